$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.957.23"
$ws.Range("E2").Value = "  +0.72%  "
$ws.Range("D3").Value = "2.318.86"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "112.46"
$ws.Range("E5").Value = "  +17.97%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "270.95"
$ws.Range("E6").Value = "  +1.59%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.628"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +0.21%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.626"
$ws.Range("E9").Value = "  +3.08%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "47.29"
$ws.Range("E10").Value = "  +7.19%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0945"
$ws.Range("E11").Value = "  +1.20%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "8.89"
$ws.Range("E12").Value = "  +15.12%  "
$ws.Range("E13").Value = "  +2.22%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "15.84"
$ws.Range("E14").Value = "  +4.40%  "
$ws.Range("D15").Value = "2.665.38"
$ws.Range("E15").Value = "  +1.30%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.866"
$ws.Range("E16").Value = "  +1.96%  "
$ws.Range("D17").Value = "2.324.74"
$ws.Range("E17").Value = "  +1.18%  "
$ws.Range("D18").Value = "43.915.84"
$ws.Range("E18").Value = "  +0.66%  "
$ws.Range("E19").Value = "  +2.76%  "
$ws.Range("E20").Value = "  +9.49%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "72.61"
$ws.Range("E21").Value = "  +0.51%  "
$ws.Range("E22").Value = "  +6.84%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "234.76"
$ws.Range("E23").Value = "  -0.19%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "9.49"
$ws.Range("E24").Value = "  +6.31%  "
$ws.Range("E25").Value = "  +15.77%  "
$ws.Range("E26").Value = "  +0.06%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.57"
$ws.Range("E27").Value = "  +2.48%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "42.83"
$ws.Range("E28").Value = "  +13.31%  "
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("E30").Value = "  +0.38%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "178.04"
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "21.93"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0934"
$ws.Range("E33").Value = "  +5.78%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.63"
$ws.Range("E34").Value = "  +5.07%  "
$ws.Range("E35").Value = "  +1.40%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.80"
$ws.Range("E36").Value = "  +8.90%  "
$ws.Range("E37").Value = "  +4.03%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.95"
$ws.Range("E38").Value = "  +21.45%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.0360"
$ws.Range("E39").Value = "  +1.99%  "
$ws.Range("E40").Value = "  +4.79%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.40"
$ws.Range("E41").Value = "  +1.48%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "69.67"
$ws.Range("E42").Value = "  +12.20%  "
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("B44").Value = "Celestia"
$ws.Range("C44").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "12.74"
$ws.Range("E44").Value = "  +7.67%  "
$ws.Range("B45").Value = "ARBITRUM"
$ws.Range("C45").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.41"
$ws.Range("E45").Value = "  +4.42%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "5.77"
$ws.Range("E46").Value = "  +10.60%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.84"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("E48").Value = "  -0.56%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "100.28"
$ws.Range("E49").Value = "  +2.30%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.22"
$ws.Range("E50").Value = "  +3.05%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.465"
$ws.Range("E51").Value = "  +10.78%  "
